$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Client / recipient block (rows 13-15) ---
$ws.Range("A13").Value = "Annaliza Mendoza Banawa"
$ws.Range("A14").Value = "ACCESS WATER INTEGRATORS AND EQUIPMENT PHILS. INC."
$ws.Range("A15").Value = "`nUnit B9 2nd Floor Regalena Bldg. National Highway Brgy. Turbina`nCalamba City`nLaguna"

# B43 originally mirrored the same company name text as A14 (shared text) - keep them in sync
$ws.Range("B43").Value = "ACCESS WATER INTEGRATORS AND EQUIPMENT PHILS. INC."

# --- ABC amount ---
$ws.Range("D23").Value = "PHP65,500.00"

# --- Purpose (new text added under the Purpose heading) ---
$ws.Range("D24").Value = "LOREM IPSUM"

# --- End user ---
$ws.Range("D28").Value = "ORD"

# --- Turn on "Wrap Text" for column A (mirrors the new wrap-enabled style) ---
$ws.Range("A1").WrapText = $true

# --- Restore the cursor/selection to A1 ---
[void]$ws.Range("A1").Select()
